$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (shifted one day forward / recalculated forecast), rows 2-15
$data = @{
    2  = @(46021, 10566.2687222002, 10821.6763103476, 19152.26, 7859.114829516, -19.6445358390159)
    3  = @(46022, 11016.6778995919, 11483.1171731316, 11232.26, 8642.43261808591, 370.55374130073)
    4  = @(46023, 5165.20006366834, 9204.49815524882, 12075.86, 8787.40063138953, 246.501616109931)
    5  = @(46024, 11651.2543252931, 11033.9950147314, 12075.86, 9010.83963541851, 332.040610422911)
    6  = @(46025, 5001.72215491993, 7795.63697117613, 12075.86, 8293.65182536117, 167.226199855721)
    7  = @(46026, 4936.28414567818, 7588.32979127335, 12075.86, 8450.06860884805, 165.105766671725)
    8  = @(46027, 13228.8164390321, 12137.103661072, 12075.86, 9108.313842858, 382.064895997083)
    9  = @(46028, 5295.53439346162, 8871.33861691519, 12075.86, 8688.196528562, 228.486464394883)
    10 = @(46029, 13228.8164390321, 13291.5461891951, 12075.86, 9108.313842858, 430.166668002213)
    11 = @(46030, 13228.8164390321, 13060.4285202311, 12075.86, 9108.313842858, 420.536765128711)
    12 = @(46031, 13228.8164390321, 12016.087873175, 12075.86, 9108.313842858, 377.022571501375)
    13 = @(46032, 5461.34628757431, 8542.74302157848, 12075.86, 8700.47192409367, 215.306456069673)
    14 = @(46033, 5295.53439346162, 8450.4611756482, 12075.86, 8688.196528562, 210.949904342092)
    15 = @(46034, 13381.7264569213, 12308.2854602627, 12075.86, 9107.02556866384, 389.143792871941)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("A$row").Value = $values[0]
    $ws.Range("B$row").Value = $values[1]
    $ws.Range("C$row").Value = $values[2]
    $ws.Range("D$row").Value = $values[3]
    $ws.Range("E$row").Value = $values[4]
    $ws.Range("F$row").Value = $values[5]
}
